$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Attendance table "טבלה3": column L = "Sun 1/2", column M = "Thu 5/2".
# Mark attendance (1) for students who attended the new "Thu 5/2" session
# (and, for a couple of students, the "Sun 1/2" session they'd missed).
$ws.Range("M5").Value = 1
$ws.Range("M6").Value = 1
$ws.Range("M7").Value = 1
$ws.Range("M12").Value = 1
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1
$ws.Range("M14").Value = 1
$ws.Range("M15").Value = 1
$ws.Range("M16").Value = 1
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1

# Clear stray explicit 0s left over in the "Thu 29/1"/"Thu 22/1" etc. columns
# so the cells go back to blank (keeping their existing cell formatting).
$ws.Range("J8").ClearContents()
$ws.Range("J9").ClearContents()
$ws.Range("G10").ClearContents()
$ws.Range("H10").ClearContents()
$ws.Range("I10").ClearContents()
$ws.Range("J10").ClearContents()
$ws.Range("K10").ClearContents()
$ws.Range("I11").ClearContents()
$ws.Range("J11").ClearContents()
$ws.Range("K11").ClearContents()
$ws.Range("K12").ClearContents()
$ws.Range("K16").ClearContents()

# Leave the cursor where the author left it.
$ws.Range("L12").Select()
